$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 635. This pushes the existing rows 635..681
# down to 636..682, preserving all their data and formatting (matching the
# dimension change from A1:R681 to A1:R682 in the target diff).
$ws.Rows("635").Insert()

# Populate the newly inserted row 635 with the new weekly record.
$ws.Cells.Item(635, 1).Value  = 3
$ws.Cells.Item(635, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(635, 3).Value  = "Coquimbo"
$ws.Cells.Item(635, 4).Value  = 45265
$ws.Cells.Item(635, 5).Value  = 5
$ws.Cells.Item(635, 6).Value  = 100112043
$ws.Cells.Item(635, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(635, 8).Value  = "Sin especificar"
$ws.Cells.Item(635, 9).Value  = "Primera"
$ws.Cells.Item(635, 10).Value = 50
$ws.Cells.Item(635, 11).Value = 19000
$ws.Cells.Item(635, 12).Value = 19000
$ws.Cells.Item(635, 13).Value = 19000
$ws.Cells.Item(635, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(635, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(635, 16).Value = 317
$ws.Cells.Item(635, 17).Value = 60
$ws.Cells.Item(635, 18).Value = "Hortaliza"
